# LV_T1516_GiftLog_ClientGiftPreApprovalPageRecipientsExceedsYearlyGiftAllowance.xlsx
# Merge - GiftLogs - Tcs. Pages & Test Data - 9th Oct 2025

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Content change: rename test user "Melissa Zatta" -> "Julie Carthane"
#    This shared string is used both on the Users sheet (A2) and on the
#    GiftLog sheet (B2, the "SubmittedFor" sample value), so update both
#    occurrences so the workbook keeps a single shared-string entry and
#    mark those cells as wrapped text (matches the style applied in Excel).
# ---------------------------------------------------------------------------
$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Range("A2").Value = "Julie Carthane"
$wsUsers.Range("A2").WrapText = $true

$wsGiftLog = $wb.Worksheets.Item("GiftLog")
$wsGiftLog.Range("B2").Value = "Julie Carthane"
$wsGiftLog.Range("B2").WrapText = $true

# ---------------------------------------------------------------------------
# 2) Navigate around the workbook the same way the author apparently did,
#    updating the stored selection on each sheet and leaving "Users" as the
#    final active sheet/tab.
# ---------------------------------------------------------------------------
$wsAppName = $wb.Worksheets.Item("AppName")
$wsAppName.Activate()
[void]$wsAppName.Range("E8").Select()

$wsModuleName = $wb.Worksheets.Item("ModuleName")
$wsModuleName.Activate()
[void]$wsModuleName.Range("D5").Select()

$wsGiftLog.Activate()
[void]$wsGiftLog.Range("D6").Select()

$wsContact = $wb.Worksheets.Item("Contact")
$wsContact.Activate()
[void]$wsContact.Range("A2").Select()

$wsContactTypes = $wb.Worksheets.Item("ContactTypes")
$wsContactTypes.Activate()
[void]$wsContactTypes.Range("B2").Select()

$wsGiftLogCurrency = $wb.Worksheets.Item("GiftLog_Currency")
$wsGiftLogCurrency.Activate()
[void]$wsGiftLogCurrency.Range("L21").Select()

$wsUsers.Activate()
[void]$wsUsers.Range("B20").Select()
